# cantor.pptx - tweak TP_preserved_under_isomorphism / typos FP_partial_order_short_answer
#
# Slide 28 ("Real Numbers Uncountable") has a text box (shape "TextBox 5")
# whose second paragraph reads:  "3 1/3  = 111.010101..."
# and third paragraph reads:     "maps to 111010101..."
#
# The fix:
#   "3 1/3"        -> "7 1/3"         (typo: should start with 7, not 3)
#   "111.010101..."-> split run so the trailing ellipsis is its own run
#   "111010101..." -> split run so the trailing ellipsis is its own run
#
# (PowerPoint naturally breaks a run into several <a:r> runs at each edited
#  character boundary, which is what we reproduce below via TextRange.Characters.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(28)
$shp = $s.Shapes.Item(3)
$tr = $shp.TextFrame.TextRange

$ellipsis = [char]0x2026

# "3 1/3" -> "7" / " " / "1/3"
$tr.Characters(22, 1).Text = "7"
$tr.Characters(24, 3).Text = "1/3"

# "111.010101…" -> "111.010101" / "…"
$tr.Characters(41, 1).Text = $ellipsis

# "111010101…" -> "111010101" / "…"
$tr.Characters(60, 1).Text = $ellipsis
